# FEINmoreThan9Error.xlsx — "Added RAD Test Cases and data for MRF."
#
# The RAD (test-automation) run timestamps recorded in column B (rows 2-13)
# are refreshed to reflect the latest Katalon test execution
# (Mon Oct 02 16:26:47 EDT 2023 .. Mon Oct 02 16:28:43 EDT 2023), replacing
# the previous run's timestamps. No other cell content, formatting, or
# layout changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = "Mon Oct 02 16:26:47 EDT 2023"
$ws.Range("B3").Value  = "Mon Oct 02 16:26:57 EDT 2023"
$ws.Range("B4").Value  = "Mon Oct 02 16:27:08 EDT 2023"
$ws.Range("B5").Value  = "Mon Oct 02 16:27:18 EDT 2023"
$ws.Range("B6").Value  = "Mon Oct 02 16:27:29 EDT 2023"
$ws.Range("B7").Value  = "Mon Oct 02 16:27:40 EDT 2023"
$ws.Range("B8").Value  = "Mon Oct 02 16:27:50 EDT 2023"
$ws.Range("B9").Value  = "Mon Oct 02 16:28:01 EDT 2023"
$ws.Range("B10").Value = "Mon Oct 02 16:28:12 EDT 2023"
$ws.Range("B11").Value = "Mon Oct 02 16:28:22 EDT 2023"
$ws.Range("B12").Value = "Mon Oct 02 16:28:33 EDT 2023"
$ws.Range("B13").Value = "Mon Oct 02 16:28:43 EDT 2023"
